$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new day's data as row 12 first (while A11 still carries the
# special "last row" date style), so we can copy that style onto A12.
$ws.Range("A12").Value = 45961
$ws.Range("B12").Value = 23
$ws.Range("C12").Value = 32
$ws.Range("D12").Value = 25

# A12 becomes the new "last row" - give it the style A11 used to have.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A12").Value = 45961

# Row 11 is no longer the last row, so restyle A11 like the other regular
# date cells above it (e.g. A10).
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A11").Value = 45960
